# Apply the "property boat&car done" edit to the 汽車 (Car) sheet (sheet3).
#
# Summary of the change:
#  - The car sheet previously had no header row labels (B1/D1/E1/F1 held
#    stray string values, C1/G1 held raw numbers) and only 6 data columns
#    (B..G) for rows 2-3.
#  - The sheet is extended to 14 columns (B..N) with a proper header row
#    matching the other property sheets (name, capacity, owner,
#    register_date, register_reason, acquire_value, property_category,
#    category, date, legislator_name, legislator_id, source_file, index),
#    and each data row is populated with the property_category/category/
#    date/legislator_name/legislator_id/source_file/index metadata columns
#    (H..N) that the other sheets already carry.
#  - Row 2's name changes from the "BluebirdSylphyGllTS" to reflect the
#    correct register_date value ("98年10月12日") and row 3's
#    register_date becomes "99年05月03日".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Header row (row 1) ----
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

# The "date" column (J) holds the literal text "2012-04-06". Force a text
# number format first so Excel doesn't auto-convert the ISO-looking string
# into a real date serial value.
$ws.Cells.Item(2,10).NumberFormat = "@"
$ws.Cells.Item(3,10).NumberFormat = "@"

# ---- Row 2 (car #30, BluebirdSylphyGllTS) ----
$ws.Cells.Item(2,2).Value = "BluebirdSylphyGllTS(客車）"
$ws.Cells.Item(2,3).Value = 1997
$ws.Cells.Item(2,4).Value = "吳美惠"
$ws.Cells.Item(2,5).Value = "98年10月12日"
$ws.Cells.Item(2,6).Value = "買賣"
$ws.Cells.Item(2,7).Value = 740000
$ws.Cells.Item(2,8).Value = "land"
$ws.Cells.Item(2,9).Value = "normal"
$ws.Cells.Item(2,10).Value = "2012-04-06"
$ws.Cells.Item(2,11).Value = "林國正"
$ws.Cells.Item(2,12).Value = 1742
$ws.Cells.Item(2,13).Value = "tmpd6491"
$ws.Cells.Item(2,14).Value = 30

# ---- Row 3 (car #31, 納智捷L91ML) ----
$ws.Cells.Item(3,2).Value = "納智捷L91ML(客車）"
$ws.Cells.Item(3,3).Value = 2198
$ws.Cells.Item(3,4).Value = "林國正"
$ws.Cells.Item(3,5).Value = "99年05月03日"
$ws.Cells.Item(3,6).Value = "買賣"
$ws.Cells.Item(3,7).Value = 800000
$ws.Cells.Item(3,8).Value = "land"
$ws.Cells.Item(3,9).Value = "normal"
$ws.Cells.Item(3,10).Value = "2012-04-06"
$ws.Cells.Item(3,11).Value = "林國正"
$ws.Cells.Item(3,12).Value = 1742
$ws.Cells.Item(3,13).Value = "tmpd6491"
$ws.Cells.Item(3,14).Value = 31
